$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate the most recent week block (rows 86-98) into new rows 101-113
#    for the new week being logged. Paste formats first (preserves the exact
#    style indices used by the source block) then paste values separately.
# ---------------------------------------------------------------------------
$srcBlock = $ws.Range("A86:I98")
$dstBlock = $ws.Range("A101:I113")

$srcBlock.Copy()
$dstBlock.PasteSpecial(-4122)   # xlPasteFormats
$srcBlock.Copy()
$dstBlock.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# Re-create the merged header cell for the new block (B101:H101), then
# re-apply the source header formatting so the merge doesn't regenerate new
# (border-split) styles.
$ws.Range("B101:H101").Merge()
$ws.Range("B86:H86").Copy()
$ws.Range("B101:H101").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Update the totals/hours that changed within the existing week block
#    (rows 86-98).
# ---------------------------------------------------------------------------
$ws.Range("I89").Value2 = 1
$ws.Range("F91").Value2 = 1
$ws.Range("I92").Value2 = 1
$ws.Range("I93").Value2 = 1
$ws.Range("I94").Value2 = 2
$ws.Range("F95").Value2 = 0.5
$ws.Range("G95").Value2 = 1
$ws.Range("I95").Value2 = 5.5
$ws.Range("C98").Value2 = 3.5
$ws.Range("D98").Value2 = 2
$ws.Range("E98").Value2 = 5
$ws.Range("F98").Value2 = 1.5
$ws.Range("G98").Value2 = 1
$ws.Range("I98").Value2 = 13

# ---------------------------------------------------------------------------
# 3. Fill in the new week's actual data on top of the duplicated block
#    (rows 101-113): dates and hours logged per task.
# ---------------------------------------------------------------------------
$ws.Range("B102").Value2 = 45601
$ws.Range("C102").Value2 = 45602
$ws.Range("D102").Value2 = 45603
$ws.Range("E102").Value2 = 45238
$ws.Range("F102").Value2 = 45239
$ws.Range("G102").Value2 = 45240
$ws.Range("H102").Value2 = 45241

$ws.Range("D104").Value2 = 1
$ws.Range("I104").Value2 = 2

$ws.Range("C105").Value2 = 1
$ws.Range("I105").Value2 = 1

$ws.Range("I106").Value2 = 0

$ws.Range("I107").Value2 = 1

$ws.Range("I108").Value2 = 1

$ws.Range("C109").ClearContents()
$ws.Range("I109").Value2 = 1

$ws.Range("F110").Value2 = 1
$ws.Range("G110").Value2 = 6
$ws.Range("I110").Value2 = 11

$ws.Range("C113").Value2 = 2
$ws.Range("D113").Value2 = 3
$ws.Range("E113").Value2 = 4
$ws.Range("F113").Value2 = 1
$ws.Range("G113").Value2 = 6
$ws.Range("I113").Value2 = 16

# ---------------------------------------------------------------------------
# 4. Scroll the view down to the newly added block and select the last cell.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
[void]$ws.Range("I113").Select()
